$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 3820017794
$ws.Range("C3").Value = 30
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "70009643"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "20220514"

# Row 4
$ws.Range("B4").Value = 3820017685
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "70006399"
$ws.Range("F4").Value = "ZRPL"
$ws.Range("G4").Value = 418
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "20220219"

# Row 5
$ws.Range("B5").Value = 3820017673
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "70006401"
$ws.Range("F5").Value = "ZRPL"
$ws.Range("G5").Value = 1023
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "20220129"

# Row 6
$ws.Range("B6").Value = 3820017810
$ws.Range("C6").Value = 50
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "70007719"
$ws.Range("F6").Value = "ZRPL"
$ws.Range("G6").Value = 50
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "20220122"

# Row 7
$ws.Range("B7").Value = 3820017810
$ws.Range("C7").Value = 40
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70006908"
$ws.Range("G7").Value = 75
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "20220122"

# Row 8
$ws.Range("B8").Value = 3820017687
$ws.Range("C8").Value = 550
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "70006417"
$ws.Range("F8").Value = "ZRRL"
$ws.Range("G8").Value = 600
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "20220212"

# Row 9
$ws.Range("B9").Value = 3820017687
$ws.Range("C9").Value = 550
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "70002173"
$ws.Range("F9").Value = "ZRER"
$ws.Range("G9").Value = 6
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "20220212"

# Row 10
$ws.Range("B10").Value = 3820017687
$ws.Range("C10").Value = 490
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "70007003"
$ws.Range("F10").Value = "ZRRL"
$ws.Range("G10").Value = 498
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "20220212"

# Row 11
$ws.Range("B11").Value = 3820017687
$ws.Range("C11").Value = 490
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "70002172"
$ws.Range("F11").Value = "ZRER"
$ws.Range("G11").Value = 6
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "20220212"

# Row 12
$ws.Range("B12").Value = 3820017687
$ws.Range("C12").Value = 430
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "70007003"
$ws.Range("F12").Value = "ZRRL"
$ws.Range("G12").Value = 250
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "20220212"

# Row 13
$ws.Range("B13").Value = 3820017687
$ws.Range("C13").Value = 380
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "70006920"
$ws.Range("F13").Value = "ZRPL"
$ws.Range("G13").Value = 1000
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "20220212"

# Row 14
$ws.Range("B14").Value = 3820017687
$ws.Range("C14").Value = 310
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "70006581"
$ws.Range("G14").Value = 300
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "20220212"

# Row 15
$ws.Range("B15").Value = 3820017687
$ws.Range("C15").Value = 250
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "70006402"
$ws.Range("F15").Value = "ZRPL"
$ws.Range("G15").Value = 750
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "20220212"

# Row 16
$ws.Range("B16").Value = 3820017687
$ws.Range("C16").Value = 190
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70006401"
$ws.Range("F16").Value = "ZRPL"
$ws.Range("G16").Value = 930
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "20220212"

# Row 17
$ws.Range("B17").Value = 3820017687
$ws.Range("C17").Value = 160
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70006244"
$ws.Range("F17").Value = "ZRPL"
$ws.Range("G17").Value = 200
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "20220212"

# Row 18
$ws.Range("B18").Value = 3820017687
$ws.Range("C18").Value = 100
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70007901"
$ws.Range("F18").Value = "ZRRL"
$ws.Range("G18").Value = 42
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "20220212"

# Row 19
$ws.Range("B19").Value = 3820017687
$ws.Range("C19").Value = 100
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70001855"
$ws.Range("F19").Value = "ZRER"
$ws.Range("G19").Value = 6
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = "20220212"

# Row 20
$ws.Range("B20").Value = 3820017687
$ws.Range("C20").Value = 40
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70006593"
$ws.Range("F20").Value = "ZRRL"
$ws.Range("G20").Value = 78
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "20220212"

# Row 21
$ws.Range("B21").Value = 3820017687
$ws.Range("C21").Value = 40
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70002433"
$ws.Range("F21").Value = "ZRER"
$ws.Range("G21").Value = 6
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "20220212"

# Row 22
$ws.Range("B22").Value = 3820017687
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70007904"
$ws.Range("F22").Value = "ZRRL"
$ws.Range("G22").Value = 132
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = "20220212"

# Row 23
$ws.Range("B23").Value = 3820017687
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70003199"
$ws.Range("F23").Value = "ZRER"
$ws.Range("G23").Value = 6
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = "20220212"

# Row 24
$ws.Range("B24").Value = 3820017809
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70008430"
$ws.Range("F24").Value = "ZRRL"
$ws.Range("G24").Value = 50
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = "20220212"

# Row 25
$ws.Range("B25").Value = 3820017809
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70004965"
$ws.Range("F25").Value = "ZRER"
$ws.Range("G25").Value = 2
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value = "20220212"

# Row 26
$ws.Range("B26").Value = 3820017684
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70006402"
$ws.Range("F26").Value = "ZRPL"
$ws.Range("G26").Value = 525
$ws.Range("J26").NumberFormat = "@"
$ws.Range("J26").Value = "20220122"

# Row 27
$ws.Range("B27").Value = 3820017729
$ws.Range("C27").Value = 10
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70005393"
$ws.Range("F27").Value = "ZRFC"
$ws.Range("G27").Value = 3261
$ws.Range("J27").NumberFormat = "@"
$ws.Range("J27").Value = "20220122"

# Row 28
$ws.Range("B28").Value = 3820017729
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "70004548"
$ws.Range("F28").Value = "ZRFC"
$ws.Range("G28").Value = 3261
$ws.Range("J28").NumberFormat = "@"
$ws.Range("J28").Value = "20220205"

# Row 29
$ws.Range("B29").Value = 3820017729
$ws.Range("C29").Value = 10
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "70004550"
$ws.Range("F29").Value = "ZREA"
$ws.Range("G29").Value = 3261
$ws.Range("J29").NumberFormat = "@"
$ws.Range("J29").Value = "20220326"

# Row 30: data removed (item no longer present), clear the row back to blank template state
$ws.Range("A30:N30").ClearContents()
